# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets
# to reflect newly refreshed counts, matching the upstream gh-pages data
# refresh at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 90
$ws1.Range("F3").Value = 4037
$ws1.Range("F4").Value = 2365
$ws1.Range("F8").Value = 26
$ws1.Range("F11").Value = 75
$ws1.Range("F12").Value = 132
$ws1.Range("F13").Value = 1507
$ws1.Range("F15").Value = 2858

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 90
$ws4.Range("F3").Value = 4037
$ws4.Range("F4").Value = 2365
$ws4.Range("F8").Value = 26
$ws4.Range("F12").Value = 75
$ws4.Range("F13").Value = 132
$ws4.Range("F16").Value = 1507
$ws4.Range("F18").Value = 2858
